# Auto commit at 2025-09-16  9:39:29.59
# Updates Metrics figures for the new day, fills in the "today" deltas
# (B3:B6 on the today sheet), and moves the active-tab/selection state
# from the "today" sheet back to "Metrics".

$wb = $excel.ActiveWorkbook

$wsToday   = $wb.Worksheets.Item("today")
$wsMetrics = $wb.Worksheets.Item("Metrics")

# --- Metrics sheet: refreshed metric values (B2:B13) ---
$wsMetrics.Range("B2").Value  = 238464
$wsMetrics.Range("B3").Value  = 193014.88
$wsMetrics.Range("B4").Value  = 75320.39
$wsMetrics.Range("B5").Value  = 9291
$wsMetrics.Range("B6").Value  = 4157714.88
$wsMetrics.Range("B7").Value  = 3520542.3599999994
$wsMetrics.Range("B8").Value  = 1204686.07
$wsMetrics.Range("B9").Value  = 160451
$wsMetrics.Range("B10").Value = 32623038.680999827
$wsMetrics.Range("B11").Value = 19550412.43
$wsMetrics.Range("B12").Value = 11486394.960000001
$wsMetrics.Range("B13").Value = 1258078

# --- today sheet: newly-filled daily delta cells (B3:B6) ---
$wsToday.Range("B3").Value = 15925.96
$wsToday.Range("B4").Value = 12884.98
$wsToday.Range("B5").Value = 5115.25
$wsToday.Range("B6").Value = 621

# --- View/selection state ---
# While "today" is still the active sheet, record its new selection.
[void]$wsToday.Range("F11:F22").Select()

# Switch the active sheet back to "Metrics" and set its selection.
$wsMetrics.Activate()
[void]$wsMetrics.Range("E10").Select()
